$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 41669456
$ws.Range("I100").Value = 83335040
$ws.Range("J100").Value = 3875
$ws.Range("K100").Value = 83335040
$ws.Range("L100").Value = 3875
$ws.Range("M100").Value = -83334499
$ws.Range("N100").Value = -4957

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 599.5769
$ws.Range("I135").Value = 435.68182
$ws.Range("J135").Value = 1501
$ws.Range("K135").Value = 3921.13638
$ws.Range("L135").Value = 13509
$ws.Range("M135").Value = -1386.13638
$ws.Range("N135").Value = -18579

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1731.375
$ws.Range("I137").Value = 1134.7894
$ws.Range("J137").Value = 3998.4
$ws.Range("K137").Value = 3404.3682
$ws.Range("L137").Value = 11995.2
$ws.Range("M137").Value = -854.3681999999999
$ws.Range("N137").Value = -17095.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6719.1025
$ws.Range("I141").Value = 1328.6538
$ws.Range("K141").Value = 3985.9614
$ws.Range("M141").Value = 1194.0386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2851.5
$ws.Range("I61").Value = 7206
$ws.Range("J61").Value = 1400
$ws.Range("K61").Value = 7206
$ws.Range("L61").Value = 1400
$ws.Range("M61").Value = -6994
$ws.Range("N61").Value = -1824

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I74").Value = 2079.8462
$ws.Range("J74").Value = 1000
$ws.Range("K74").Value = 2079.8462
$ws.Range("L74").Value = 1000
$ws.Range("M74").Value = -1205.8462
$ws.Range("N74").Value = -2748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I77").Value = 2079.8462
$ws.Range("J77").Value = 1000
$ws.Range("K77").Value = 10399.231
$ws.Range("L77").Value = 5000
$ws.Range("M77").Value = -6031.231
$ws.Range("N77").Value = -13736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3574.5715
$ws.Range("I122").Value = 2341
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 7023
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -4573
$ws.Range("N122").Value = -18399.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2851.5
$ws.Range("I136").Value = 7206
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 21618
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -19068
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 68730
$ws.Range("J139").Value = 68730
$ws.Range("L139").Value = 68730
$ws.Range("N139").Value = -79010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1607.7727
$ws.Range("I99").Value = 1143.0769
$ws.Range("J99").Value = 2279
$ws.Range("K99").Value = 1143.0769
$ws.Range("L99").Value = 2279
$ws.Range("M99").Value = 354.9231
$ws.Range("N99").Value = -5275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 56914.2
$ws.Range("J28").Value = 56914.2
$ws.Range("L28").Value = 56914.2
$ws.Range("N28").Value = -57404.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2455.2222
$ws.Range("I58").Value = 1910
$ws.Range("J58").Value = 3136.75
$ws.Range("K58").Value = 1910
$ws.Range("L58").Value = 3136.75
$ws.Range("M58").Value = -1707
$ws.Range("N58").Value = -3542.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2430.963
$ws.Range("I134").Value = 2379.6
$ws.Range("K134").Value = 7138.799999999999
$ws.Range("M134").Value = -4603.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2455.2222
$ws.Range("I136").Value = 1910
$ws.Range("J136").Value = 3136.75
$ws.Range("K136").Value = 5730
$ws.Range("L136").Value = 9410.25
$ws.Range("M136").Value = -3180
$ws.Range("N136").Value = -14510.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 971.8570999999999
$ws.Range("I92").Value = 926
$ws.Range("J92").Value = 1033
$ws.Range("K92").Value = 2778
$ws.Range("L92").Value = 3099
$ws.Range("M92").Value = -1530
$ws.Range("N92").Value = -5595

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 885.7143
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 885.7143
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2657.1429
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3649.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1249.25
$ws.Range("I98").Value = 2503
$ws.Range("J98").Value = 831.3333
$ws.Range("K98").Value = 7509
$ws.Range("L98").Value = 2493.9999
$ws.Range("M98").Value = -6011
$ws.Range("N98").Value = -5489.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4277.1875
$ws.Range("J131").Value = 6354.7617
$ws.Range("L131").Value = 19064.2851
$ws.Range("N131").Value = -29144.2851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1833.3334
$ws.Range("I113").Value = 1833.3334
$ws.Range("K113").Value = 1833.3334
$ws.Range("M113").Value = 336.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 42000
$ws.Range("J56").Value = 42000
$ws.Range("L56").Value = 42000
$ws.Range("N56").Value = -43382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1538.5
$ws.Range("I61").Value = 1601.3334
$ws.Range("J61").Value = 1500.8
$ws.Range("K61").Value = 1601.3334
$ws.Range("L61").Value = 1500.8
$ws.Range("M61").Value = -1399.3334
$ws.Range("N61").Value = -1904.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1538.5
$ws.Range("I113").Value = 1601.3334
$ws.Range("J113").Value = 1500.8
$ws.Range("K113").Value = 1601.3334
$ws.Range("L113").Value = 1500.8
$ws.Range("M113").Value = 568.6666
$ws.Range("N113").Value = -5840.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13892322
$ws.Range("I136").Value = 3944.1875
$ws.Range("J136").Value = 41669080
$ws.Range("K136").Value = 11832.5625
$ws.Range("L136").Value = 125007240
$ws.Range("M136").Value = -9282.5625
$ws.Range("N136").Value = -125012340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1643.3721
$ws.Range("I132").Value = 1444.7715
$ws.Range("J132").Value = 2512.25
$ws.Range("K132").Value = 4334.3145
$ws.Range("L132").Value = 7536.75
$ws.Range("M132").Value = -1804.3145
$ws.Range("N132").Value = -12596.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4486.3125
$ws.Range("I136").Value = 1975.6
$ws.Range("J136").Value = 5627.5454
$ws.Range("K136").Value = 5926.799999999999
$ws.Range("L136").Value = 16882.6362
$ws.Range("M136").Value = -3376.799999999999
$ws.Range("N136").Value = -21982.6362
